$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1749.5
$ws.Range("I19").Value = 1999.5
$ws.Range("J19").Value = 1499.5
$ws.Range("K19").Value = 1999.5
$ws.Range("L19").Value = 1499.5
$ws.Range("M19").Value = -1824.5
$ws.Range("N19").Value = -1849.5
$ws.Range("H41").Value = 6583.5
$ws.Range("I41").Value = 3900.6
$ws.Range("K41").Value = 3900.6
$ws.Range("M41").Value = -3460.6
$ws.Range("H55").Value = 390.8
$ws.Range("J55").Value = 538
$ws.Range("L55").Value = 538
$ws.Range("N55").Value = -966
$ws.Range("H86").Value = 2492.7144
$ws.Range("J86").Value = 2949.2
$ws.Range("L86").Value = 2949.2
$ws.Range("N86").Value = -5195.2
$ws.Range("H89").Value = 2492.7144
$ws.Range("J89").Value = 2949.2
$ws.Range("L89").Value = 14746
$ws.Range("N89").Value = -25978
$ws.Range("H106").Value = 3998.3333
$ws.Range("I106").Value = 3998.3333
$ws.Range("K106").Value = 3998.3333
$ws.Range("M106").Value = -3367.3333
$ws.Range("H111").Value = 666
$ws.Range("I111").Value = 666
$ws.Range("K111").Value = 1998
$ws.Range("M111").Value = 1069
$ws.Range("H112").Value = 3049.5
$ws.Range("J112").Value = 4999
$ws.Range("L112").Value = 14997
$ws.Range("N112").Value = -17213
$ws.Range("H116").Value = 5333.3335
$ws.Range("I116").Value = 3500
$ws.Range("K116").Value = 3500
$ws.Range("M116").Value = -58
$ws.Range("H135").Value = 1589.8572
$ws.Range("I135").Value = 1589.8572
$ws.Range("K135").Value = 14308.7148
$ws.Range("M135").Value = -11773.7148
$ws.Range("H138").Value = 3999
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2867.75
$ws.Range("I61").Value = 2848.923
$ws.Range("J61").Value = 2949.3333
$ws.Range("K61").Value = 2848.923
$ws.Range("L61").Value = 2949.3333
$ws.Range("M61").Value = -2636.923
$ws.Range("N61").Value = -3373.3333
$ws.Range("H74").Value = 10500
$ws.Range("I74").Value = 6000
$ws.Range("K74").Value = 6000
$ws.Range("M74").Value = -5126
$ws.Range("H77").Value = 10500
$ws.Range("I77").Value = 6000
$ws.Range("K77").Value = 30000
$ws.Range("M77").Value = -25632
$ws.Range("H97").Value = 3891.3
$ws.Range("I97").Value = 3016.1428
$ws.Range("J97").Value = 5933.3335
$ws.Range("K97").Value = 3016.1428
$ws.Range("L97").Value = 5933.3335
$ws.Range("M97").Value = -2520.1428
$ws.Range("N97").Value = -6925.3335
$ws.Range("H122").Value = 6152.773
$ws.Range("I122").Value = 5968.25
$ws.Range("J122").Value = 7998
$ws.Range("K122").Value = 17904.75
$ws.Range("L122").Value = 23994
$ws.Range("M122").Value = -15454.75
$ws.Range("N122").Value = -28894
$ws.Range("H134").Value = 79499
$ws.Range("J134").Value = 79499
$ws.Range("L134").Value = 79499
$ws.Range("N134").Value = -89639
$ws.Range("H136").Value = 2867.75
$ws.Range("I136").Value = 2848.923
$ws.Range("J136").Value = 2949.3333
$ws.Range("K136").Value = 8546.769
$ws.Range("L136").Value = 8847.999899999999
$ws.Range("M136").Value = -5996.769
$ws.Range("N136").Value = -13947.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 6625
$ws.Range("J106").Value = 6625
$ws.Range("L106").Value = 6625
$ws.Range("N106").Value = -9149
$ws.Range("H134").Value = 5449.25
$ws.Range("I134").Value = 914.1429000000001
$ws.Range("K134").Value = 2742.4287
$ws.Range("M134").Value = -207.4287000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 349.66666
$ws.Range("J7").Value = 412.25
$ws.Range("L7").Value = 412.25
$ws.Range("N7").Value = -638.25
$ws.Range("H31").Value = 7056.769
$ws.Range("I31").Value = 6873.6665
$ws.Range("K31").Value = 6873.6665
$ws.Range("M31").Value = -6578.6665
$ws.Range("H34").Value = 7056.769
$ws.Range("I34").Value = 6873.6665
$ws.Range("K34").Value = 6873.6665
$ws.Range("M34").Value = -6671.6665
$ws.Range("H58").Value = 453.66666
$ws.Range("I58").Value = 453.66666
$ws.Range("K58").Value = 453.66666
$ws.Range("M58").Value = -250.66666
$ws.Range("H107").Value = 700
$ws.Range("J107").Value = 499.5
$ws.Range("L107").Value = 499.5
$ws.Range("N107").Value = -4339.5
$ws.Range("H136").Value = 453.66666
$ws.Range("I136").Value = 453.66666
$ws.Range("K136").Value = 1360.99998
$ws.Range("M136").Value = 1189.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 18000
$ws.Range("M70").Value = -17685
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 18000
$ws.Range("M73").Value = -16908
$ws.Range("H80").Value = 5175
$ws.Range("I80").Value = 4350
$ws.Range("K80").Value = 13050
$ws.Range("M80").Value = -12114
$ws.Range("H83").Value = 5175
$ws.Range("I83").Value = 4350
$ws.Range("K83").Value = 39150
$ws.Range("M83").Value = -34470
$ws.Range("H98").Value = 2004
$ws.Range("J98").Value = 2004
$ws.Range("L98").Value = 6012
$ws.Range("N98").Value = -9008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7500000
$ws.Range("I10").Value = 7500000
$ws.Range("K10").Value = 7500000
$ws.Range("M10").Value = -7499831
$ws.Range("H70").Value = 6125.25
$ws.Range("J70").Value = 6750.5
$ws.Range("L70").Value = 6750.5
$ws.Range("N70").Value = -7290.5
$ws.Range("H73").Value = 6125.25
$ws.Range("J73").Value = 6750.5
$ws.Range("L73").Value = 6750.5
$ws.Range("N73").Value = -8622.5
$ws.Range("H107").Value = 725.3333
$ws.Range("I107").Value = 175
$ws.Range("J107").Value = 1000.5
$ws.Range("K107").Value = 175
$ws.Range("L107").Value = 1000.5
$ws.Range("M107").Value = 1745
$ws.Range("N107").Value = -4840.5
$ws.Range("H122").Value = 1100.5555
$ws.Range("I122").Value = 986.5714
$ws.Range("J122").Value = 1499.5
$ws.Range("K122").Value = 2959.7142
$ws.Range("L122").Value = 4498.5
$ws.Range("M122").Value = -509.7142000000003
$ws.Range("N122").Value = -9398.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 20005
$ws.Range("I3").Value = 20005
$ws.Range("K3").Value = 20005
$ws.Range("M3").Value = -19893
$ws.Range("H7").Value = 5799.6
$ws.Range("I7").Value = 5799.6
$ws.Range("K7").Value = 5799.6
$ws.Range("M7").Value = -5687.6
$ws.Range("H15").Value = 20005
$ws.Range("I15").Value = 20005
$ws.Range("K15").Value = 20005
$ws.Range("M15").Value = -19835
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = $null
$ws.Range("N22").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = $null
$ws.Range("N27").Value = 0
$ws.Range("H46").Value = 2563
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2563
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = $null
$ws.Range("M46").Value = 2563
$ws.Range("N46").Value = -2939
$ws.Range("H93").Value = 5032.6924
$ws.Range("I93").Value = 5032.6924
$ws.Range("K93").Value = 5032.6924
$ws.Range("M93").Value = -3784.6924
$ws.Range("H94").Value = 36876.332
$ws.Range("I94").Value = 20300
$ws.Range("J94").Value = 45164.5
$ws.Range("K94").Value = 20300
$ws.Range("L94").Value = 45164.5
$ws.Range("M94").Value = -19624
$ws.Range("N94").Value = -46516.5
$ws.Range("H122").Value = 4916.5
$ws.Range("I122").Value = 4749.75
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 14249.25
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -11799.25
$ws.Range("N122").Value = -20650
$ws.Range("H126").Value = 5799.6
$ws.Range("I126").Value = 5799.6
$ws.Range("K126").Value = 17398.8
$ws.Range("M126").Value = -14928.8
$ws.Range("H132").Value = 5259.091
$ws.Range("I132").Value = 5261.1113
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 15783.3339
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -13253.3339
$ws.Range("N132").Value = -20810
$ws.Range("H136").Value = 5414.143
$ws.Range("I136").Value = 5414.143
$ws.Range("K136").Value = 16242.429
$ws.Range("M136").Value = -13692.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 2826
$ws.Range("I14").Value = 1152
$ws.Range("K14").Value = 1152
$ws.Range("M14").Value = -984
$ws.Range("H74").Value = 42995
$ws.Range("J74").Value = 42995
$ws.Range("L74").Value = 42995
$ws.Range("N74").Value = -44867
$ws.Range("H77").Value = 42995
$ws.Range("J77").Value = 42995
$ws.Range("L77").Value = 128985
$ws.Range("N77").Value = -138345
$ws.Range("H107").Value = 5300
$ws.Range("I107").Value = 1900
$ws.Range("K107").Value = 5700
$ws.Range("M107").Value = -3780
$ws.Range("H132").Value = 4425.7334
$ws.Range("I132").Value = 3494.6667
$ws.Range("J132").Value = 8150
$ws.Range("K132").Value = 10484.0001
$ws.Range("L132").Value = 24450
$ws.Range("M132").Value = -7954.000100000001
$ws.Range("N132").Value = -29510
